$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in cell A2: "grouped victim" -> "groped victim"
$ws.Range("A2").Value = "Grade 10 student believed a black van followed as she walked to school from her resident on three occasions. This morning as the victim was walking to school. An unknown suspect jumped out from the bush and groped victim from behind, then ran away. The parents were contacted and a safety plan has been put in place including for victim to not walk to school alone in the mornings for the time being. Preliminary area search negative for the vehicle and several CCTV potential locations observed."

# Update the active selection to A5
$ws.Range("A5").Select()
